# Edit timetable and link to stage one report in group agreement
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column F (Research Skills) content rewrite for weeks 1-10 ---
$ws.Range("F3").Value  = "* Registered report and MSLQ overview"
$ws.Range("F4").Value  = "* Finding, reading, and organising journal articles`n* Group work agreement"
$ws.Range("F5").Value  = "* Introduction structure`n* Identifying the rationale"
$ws.Range("F6").Value  = "* Scientific writing`n* Paragraph structure`n* Citation placement"
$ws.Range("F7").Value  = "* Method structure`n* Researcher degrees of freedom"
$ws.Range("F9").Value  = "* Correlation results sections`n* Reporting power analyses"
$ws.Range("F10").Value = "* t-test results sections"
$ws.Range("F11").Value = "* Discussion structure"
$ws.Range("F12").Value = "* Abstract structure"

# --- Row height adjustments to match new wrapped-text content ---
$ws.Rows.Item(4).RowHeight = 48
$ws.Rows.Item(6).RowHeight = 48
$ws.Rows.Item(11).RowHeight = 16
$ws.Rows.Item(12).RowHeight = 16

# --- Update selection to match saved cursor position ---
$null = $ws.Range("F13").Select()
